# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '43.312.31'
$ws.Range('E2').Value = '  -1.00%  '
# Row 3
$ws.Range('D3').Value = '2.356.33'
$ws.Range('E3').Value = '  +5.75%  '
# Row 4
$ws.Range('E4').Value = '  -0.12%  '
# Row 5
$ws.Range('D5').Value = "'232.64"
$ws.Range('E5').Value = '  +1.27%  '
# Row 6
$ws.Range('D6').Value = "'0.650"
$ws.Range('E6').Value = '  +1.38%  '
# Row 7
$ws.Range('D7').Value = "'68.02"
$ws.Range('E7').Value = '  +7.66%  '
# Row 8
$ws.Range('E8').Value = '  -0.04%  '
# Row 9
$ws.Range('E9').Value = '  +2.54%  '
# Row 10
$ws.Range('D10').Value = "'0.0955"
$ws.Range('E10').Value = '  -0.78%  '
# Row 11
$ws.Range('E11').Value = '  +0.02%  '
# Row 12
$ws.Range('D12').Value = "'26.37"
$ws.Range('E12').Value = '  -0.77%  '
# Row 13
$ws.Range('D13').Value = '2.707.37'
$ws.Range('E13').Value = '  +5.57%  '
# Row 14
$ws.Range('E14').Value = '  -0.86%  '
# Row 15
$ws.Range('D15').Value = "'15.69"
$ws.Range('E15').Value = '  +2.23%  '
# Row 16
$ws.Range('D16').Value = "'6.24"
$ws.Range('E16').Value = '  +2.71%  '
# Row 17
$ws.Range('E17').Value = '  +1.94%  '
# Row 18
$ws.Range('D18').Value = '2.353.19'
$ws.Range('E18').Value = '  +5.55%  '
# Row 19
$ws.Range('D19').Value = '43.296.83'
$ws.Range('E19').Value = '  -0.67%  '
# Row 20
$ws.Range('D20').Value = '0.0₃0982'
$ws.Range('E20').Value = '  -0.02%  '
# Row 21
$ws.Range('D21').Value = "'74.02"
$ws.Range('E21').Value = '  +2.17%  '
# Row 22
$ws.Range('E22').Value = '  +4.25%  '
# Row 23
$ws.Range('D23').Value = "'248.34"
$ws.Range('E23').Value = '  +0.30%  '
# Row 24
$ws.Range('D24').Value = "'4.00"
$ws.Range('E24').Value = '  +17.73%  '
# Row 25
$ws.Range('E25').Value = '  +0.00%  '
# Row 26
$ws.Range('E26').Value = '  +1.77%  '
# Row 27
$ws.Range('D27').Value = "'2.22"
$ws.Range('E27').Value = '  -3.27%  '
# Row 28
$ws.Range('D28').Value = "'9.92"
$ws.Range('E28').Value = '  +0.70%  '
# Row 29
$ws.Range('D29').Value = "'22.33"
$ws.Range('E29').Value = '  +7.78%  '
# Row 30
$ws.Range('D30').Value = "'173.21"
$ws.Range('E30').Value = '  +1.65%  '
# Row 31
$ws.Range('D31').Value = "'1.52"
$ws.Range('E31').Value = '  +11.15%  '
# Row 32
$ws.Range('D32').Value = "'0.127"
$ws.Range('E32').Value = '  -6.92%  '
# Row 33
$ws.Range('E33').Value = '  +1.10%  '
# Row 34
$ws.Range('E34').Value = '  +5.99%  '
# Row 35
$ws.Range('D35').Value = "'0.0695"
$ws.Range('E35').Value = '  +0.12%  '
# Row 36
$ws.Range('D36').Value = "'5.06"
$ws.Range('E36').Value = '  +4.32%  '
# Row 37
$ws.Range('D37').Value = "'2.50"
$ws.Range('E37').Value = '  +11.05%  '
# Row 38
$ws.Range('D38').Value = "'6.50"
$ws.Range('E38').Value = '  +2.23%  '
# Row 39
$ws.Range('D39').Value = "'3.63"
$ws.Range('E39').Value = '  -0.14%  '
# Row 40
$ws.Range('D40').Value = "'0.0255"
$ws.Range('E40').Value = '  -0.79%  '
# Row 41
$ws.Range('B41').Value = 'BinanceUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D41').Value = "'1.00"
$ws.Range('E41').Value = '  -0.09%  '
# Row 42
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = "'8.94"
$ws.Range('E42').Value = '  +9.20%  '
# Row 43
$ws.Range('D43').Value = "'18.19"
$ws.Range('E43').Value = '  +6.81%  '
# Row 44
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = "'1.16"
$ws.Range('E44').Value = '  +8.96%  '
# Row 45
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = "'1.21"
$ws.Range('E45').Value = '  +2.92%  '
# Row 46
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'98.61"
$ws.Range('E46').Value = '  +2.03%  '
# Row 47
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').Value = "'4.45"
$ws.Range('E47').Value = '  +2.35%  '
# Row 48
$ws.Range('D48').Value = "'0.0950"
$ws.Range('E48').Value = '  +1.25%  '
# Row 49
$ws.Range('D49').Value = '1.444.51'
$ws.Range('E49').Value = '  +1.48%  '
# Row 50
$ws.Range('D50').Value = '2.579.98'
$ws.Range('E50').Value = '  +5.85%  '
# Row 51
$ws.Range('E51').Value = '  -2.54%  '
